# pandemic2020.xlsx - "Fix broken URLs and add death rate"
# Adds a "治愈" (recovered/cured) column right after the "死亡" (H) column
# on the "data" sheet, shifting the existing ratio-table (previously J:Q)
# one column to the right (now K:R), and nudges a couple of view/selection
# bits to match the authored session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$wsChart = $wb.Worksheets.Item("chart")

# --- Insert a new column before column I (pushes J:Q -> K:R), carrying
# --- formatting/formulas along with it (relative refs re-target correctly).
$ws.Columns("I:I").Insert()

# --- New column header + data (matches the style already on column H).
$ws.Range("I2").Value = "治愈"
$ws.Range("I3").Value = 103
$ws.Range("I4").Value = 60
$ws.Range("I5").Value = 51
$ws.Range("I6").Value = 49
$ws.Range("I7").Value = 38
$ws.Range("I8").Value = 34

# Column H/I are both narrow numeric columns - match their on-screen width.
$ws.Range("H2:I2").ColumnWidth = 4.63

# --- Selections, restored to where the author left off.
$ws.Range("L15").Select()
$wsChart.Range("J57").Select()
